# Update the cryptos price list (columns D = Price, E = Volume(1h)) for rows 2-51
# with refreshed figures, matching the GitHub Actions scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.159.59"
$ws.Range("E2").Value = "  -5.95%  "
$ws.Range("D3").Value = "'2.554.58"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'299.30"
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").Value = "'94.49"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("D7").Value = "'0.576"
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").Value = "'35.96"
$ws.Range("E10").Value = "  -6.94%  "
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").Value = "'7.76"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "'2.943.29"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "'2.539.54"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("E16").Value = "  -3.85%  "
$ws.Range("D17").Value = "'14.19"
$ws.Range("E17").Value = "  -3.77%  "
$ws.Range("D18").Value = "'43.178.38"
$ws.Range("E18").Value = "  -6.32%  "
$ws.Range("D19").Value = "'12.98"
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("D20").Value = "'0.0₃0982"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "'72.53"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'260.54"
$ws.Range("E23").Value = "  -10.11%  "
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").Value = "'29.86"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("E26").Value = "  -4.68%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = "  -6.27%  "
$ws.Range("D29").Value = "'37.53"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "'154.77"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("E36").Value = "  -4.54%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").Value = "'16.78"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("D40").Value = "'23.53"
$ws.Range("E40").Value = "  +11.32%  "
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").Value = "'0.0314"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").Value = "'3.91"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").Value = "'2.075.15"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'86.00"
$ws.Range("E46").Value = "  -9.81%  "
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").Value = "'2.800.51"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").Value = "'8.72"
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("D51").Value = "'104.55"
$ws.Range("E51").Value = "  -4.01%  "
